$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (58) with the latest month's data (01-09-2021)
# Force column A to be treated as text so the date-like string isn't
# auto-converted into a date serial number, then clear the formatting
# change so the cell keeps the default (no explicit) style.
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "01-09-2021"
$ws.Range("A58").ClearFormats()

$ws.Range("B58").Value = 5.7
$ws.Range("C58").Value = 0.4
$ws.Range("D58").Value = 6.2
